$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column cells stay plain text (several look like numbers,
# e.g. "246.80" or "40.40", and Excel would otherwise coerce them to
# floating-point numbers and drop the formatting/trailing zeros).
$priceCells = "D2","D3","D5","D6","D7","D10","D13","D14","D15","D17","D18","D20","D21","D23","D25","D28","D30","D32","D33","D37","D39","D40","D42","D43","D44","D45","D50","D51"
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "41.987.01"
$ws.Range("E2").Value = "  -1.48%  "

$ws.Range("D3").Value = "2.241.54"
$ws.Range("E3").Value = "  -1.83%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "246.80"
$ws.Range("E5").Value = "  -2.05%  "

$ws.Range("D6").Value = "0.630"
$ws.Range("E6").Value = "  +0.50%  "

$ws.Range("D7").Value = "75.21"
$ws.Range("E7").Value = "  +1.91%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  -2.81%  "

$ws.Range("D10").Value = "40.40"
$ws.Range("E10").Value = "  +2.94%  "

$ws.Range("E11").Value = "  -3.84%  "

$ws.Range("D13").Value = "0.103"
$ws.Range("E13").Value = "  -2.38%  "

$ws.Range("D14").Value = "2.577.28"
$ws.Range("E14").Value = "  -1.92%  "

$ws.Range("D15").Value = "14.83"
$ws.Range("E15").Value = "  -0.91%  "

$ws.Range("E16").Value = "  -1.60%  "

$ws.Range("D17").Value = "2.230.58"
$ws.Range("E17").Value = "  -2.22%  "

$ws.Range("D18").Value = "41.930.20"
$ws.Range("E18").Value = "  -1.47%  "

$ws.Range("E19").Value = "  -2.90%  "

$ws.Range("D20").Value = "6.14"
$ws.Range("E20").Value = "  -2.87%  "

$ws.Range("D21").Value = "71.47"
$ws.Range("E21").Value = "  -1.12%  "

$ws.Range("E22").Value = "  +2.31%  "

$ws.Range("D23").Value = "230.97"
$ws.Range("E23").Value = "  -0.80%  "

$ws.Range("E24").Value = "  +0.06%  "

$ws.Range("D25").Value = "11.30"
$ws.Range("E25").Value = "  -1.98%  "

$ws.Range("E26").Value = "  -4.97%  "

$ws.Range("E27").Value = "  -4.43%  "

$ws.Range("D28").Value = "7.25"
$ws.Range("E28").Value = "  +14.31%  "

$ws.Range("E29").Value = "  -1.63%  "

$ws.Range("D30").Value = "169.75"
$ws.Range("E30").Value = "  +1.74%  "

$ws.Range("D32").Value = "34.25"
$ws.Range("E32").Value = "  +7.15%  "

$ws.Range("D33").Value = "0.0841"
$ws.Range("E33").Value = "  +2.96%  "

$ws.Range("E34").Value = "  -5.67%  "

$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("E36").Value = "  -0.55%  "

$ws.Range("D37").Value = "4.92"
$ws.Range("E37").Value = "  +3.53%  "

$ws.Range("D39").Value = "13.58"
$ws.Range("E39").Value = "  -2.23%  "

$ws.Range("D40").Value = "5.93"
$ws.Range("E40").Value = "  -0.52%  "

$ws.Range("E41").Value = "  -6.40%  "

$ws.Range("D42").Value = "110.64"
$ws.Range("E42").Value = "  +13.10%  "

$ws.Range("D43").Value = "0.203"
$ws.Range("E43").Value = "  -5.11%  "

$ws.Range("D44").Value = "60.38"
$ws.Range("E44").Value = "  -2.51%  "

$ws.Range("D45").Value = "8.78"
$ws.Range("E45").Value = "  -4.13%  "

$ws.Range("E46").Value = "  -3.20%  "

$ws.Range("E47").Value = "  -0.33%  "

$ws.Range("E48").Value = "  -3.88%  "

$ws.Range("E49").Value = "  -1.24%  "

$ws.Range("D50").Value = "4.27"
$ws.Range("E50").Value = "  -12.03%  "

$ws.Range("D51").Value = "2.26"
$ws.Range("E51").Value = "  -0.99%  "
